$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Summary" - update summary metrics to reflect new trade #33
# -----------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.67   # Current Capital
$summary.Range("B4").Value = 0.66      # Total P&L $
$summary.Range("B5").Value = 0.4       # Total P&L %
$summary.Range("B6").Value = 33        # Total Trades
$summary.Range("B8").Value = 12        # Losing Trades
$summary.Range("B9").Value = 39.39     # Win Rate %

# -----------------------------------------------------------------
# Sheet "Strategy Status" - update MarketMaking row (row 4)
# -----------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.67     # Capital
$status.Range("D4").Value = 33         # Trades
$status.Range("E4").Value = 0.66       # P&L $
$status.Range("F4").Value = 0.67       # P&L %
$status.Range("G4").Value = 39.39      # Win Rate %

# -----------------------------------------------------------------
# Helper to append the new trade row (row 34) on a trade-log sheet
# -----------------------------------------------------------------
function Add-Trade33Row($ws) {
    $row = 34

    $ws.Cells.Item($row, 1).Value = 33

    # Date / Time columns need to stay as plain text, not get
    # auto-converted to Excel date/time serial numbers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = "12:38:33"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.26
    $ws.Cells.Item($row, 7).Value = 0.21
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -19.2308
    $ws.Cells.Item($row, 10).Value = -0.05
    $ws.Cells.Item($row, 11).Value = 100.67
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}

# -----------------------------------------------------------------
# Sheet "All Trades" - append trade #33
# -----------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade33Row $allTrades

# -----------------------------------------------------------------
# Sheet "MarketMaking" - append trade #33 (same log duplicated)
# -----------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade33Row $marketMaking
